{"js": "// Remove the standalone \"I personally examined the patient...\" paragraph\n// (the attestation/co-sign statement) that precedes the OBJECTIVE section.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst target =\n  \"I personally examined the patient separately and discussed the case \" +\n  \"with the resident/physician assistant and with any services involved \" +\n  \"in a multidisciplinary fashion. I agree with the resident/physician's \" +\n  \"assistant documentation with any exceptions noted below:\";\n\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.text.trim() === target) {\n    paragraph.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the standalone \"I personally examined the patient...\" paragraph\n# (the attestation/co-sign statement) that precedes the OBJECTIVE section.\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$found = $range.Find.Execute(\"I personally examined the patient separately\")\n\nif ($found) {\n    # Expand the found range to cover the whole paragraph (wdParagraph = 4)\n    # so the paragraph mark is removed along with its text, rather than\n    # leaving behind an empty paragraph.\n    $range.Expand(4) | Out-Null\n    $range.Delete()\n}\n"}
